$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 368, pushing existing rows 368:397 down to 369:398
$ws.Rows(368).Insert()

# Populate the newly inserted row 368 with the new weekly record
$ws.Range("A368").Value = 3
$ws.Range("B368").Value = "Femacal de La Calera"
$ws.Range("C368").Value = "Coquimbo"
$ws.Range("D368").Value = 44783
$ws.Range("E368").Value = 5
$ws.Range("F368").Value = 100112043
$ws.Range("G368").Value = "Pepino ensalada"
$ws.Range("H368").Value = "Sin especificar"
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 90
$ws.Range("K368").Value = 17000
$ws.Range("L368").Value = 18000
$ws.Range("M368").Value = 17500
$ws.Range("N368").Value = "`$/caja 70 unidades"
$ws.Range("O368").Value = "Región de Arica y Parinacota"
$ws.Range("P368").Value = 250
$ws.Range("Q368").Value = 70
$ws.Range("R368").Value = "Hortaliza"
